# Add a new "Baseline_2010-18_C205" simulation-run row to the "2010-18" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# New row 3 mirrors row 2's layout: model / run name / weather years / stats ... / weather years again.
$ws.Range("A3").Value = "CW3M"
$ws.Range("B3").Value = "Baseline_2010-18_C205"
$ws.Range("C3").Value = "2010-18"

$ws.Range("D3").Value = 677.32165200000009
$ws.Range("E3").Value = 2094.2995878888887
$ws.Range("F3").Value = 5.6902922222222223
$ws.Range("G3").Value = 190.42691511111113
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1.321501777777778
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 566.90397822222212
$ws.Range("L3").Value = 96.602825555555569
$ws.Range("M3").Value = 1624.5496012222222
$ws.Range("N3").Value = 681.18569955555552
$ws.Range("O3").Value = 15160.787543444445
$ws.Range("P3").Value = 2216.7525497777779
$ws.Range("Q3").Value = 0.18215544444444445
$ws.Range("R3").Value = -0.000045222222222222227
$ws.Range("S3").Value = "2010-18"

# Match row 2's number formats for the corresponding columns.
$ws.Range("D3:N3").NumberFormat = $ws.Range("D2:N2").NumberFormat
$ws.Range("P3").NumberFormat = $ws.Range("P2").NumberFormat
$ws.Range("Q3").NumberFormat = $ws.Range("Q2").NumberFormat
$ws.Range("R3").NumberFormat = $ws.Range("R2").NumberFormat

# O3 gets the highlighted (yellow-fill) style used elsewhere in the workbook for this metric.
$ws.Range("O3").NumberFormat = "0"
$ws.Range("O3").Interior.Color = 65535

# Widen column R (col index 18) slightly, matching the author's width tweak
# (engine quantizes ColumnWidth to 1/6-character steps; 11.6667 is the closest
# input that rounds to the saved OOXML width of ~12.44).
$ws.Columns.Item(18).ColumnWidth = 11.666666666666666

$ws.Range("B4").Select()
